# thêm trạng thái đóng mở cửa hàng
# Shift the "Thiết bị"/Device column out and append a new "Đóng cửa"/Closed
# column at the end (column X), matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 7): U7..X7 shift left by one, X7 becomes "Đóng cửa" ---
$ws.Range("U7").Value = "Đúng tuyến"
$ws.Range("V7").Value = "Chụp ảnh"
$ws.Range("W7").Value = "Đơn hàng"
$ws.Range("X7").Value = "Đóng cửa"

# --- Placeholder row (row 9): U9..X9 shift left by one, X9 becomes the new placeholder ---
$ws.Range("U9").Value = "{{ReportStoreCheckeds.SalesEmployees.Dates.Contents.Planned}}"
$ws.Range("V9").Value = "{{ReportStoreCheckeds.SalesEmployees.Dates.Contents.Image}}"
$ws.Range("W9").Value = "{{ReportStoreCheckeds.SalesEmployees.Dates.Contents.Order}}"
$ws.Range("X9").Value = "{{ReportStoreCheckeds.SalesEmployees.Dates.Contents.Closed}}"

# W9 picks up the same formatting (incl. right-hand border) that the
# previously-last column X9 already had, reusing the existing style
$ws.Range("X9").Copy()
$ws.Range("W9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the authored document
$ws.Range("E12").Select()
